$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds daily time-tracking entries in rows 2-14
# (column A = date, column B = hours worked as a time fraction).
# Append two more days of entries (rows 15 and 16), finishing the
# "cadastro de natureza de patente" log.
#
# New rows are created by duplicating the last data row (14) so the
# new cells inherit the same number formats (date / time) already
# used throughout the column, then the copied values are overwritten
# with the real data for the two new days.

$ws.Rows("14").Copy()
$ws.Rows("15").Insert(-4121)

$ws.Rows("14").Copy()
$ws.Rows("16").Insert(-4121)

# Row 15: 2013-10-21, 00:35 (0:35 -> 2.4305555555555556E-2 of a day)
$ws.Range("A15").Value = 41568
$ws.Range("B15").Value = 0.024305555555555556

# Row 16: 2013-10-22, 03:00 (3:00 -> 0.125 of a day)
$ws.Range("A16").Value = 41569
$ws.Range("B16").Value = 0.125

# Move the active selection down to the next empty cell below the table.
[void]$ws.Range("B17").Select()
